$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
